$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 12 with the new week's data
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 0.62582175925925931
$ws.Range("D12").Value = "Watched dubbed shows and movies without subs and Spanish shows with subs, also a bit of  Youtube videos about Minecraft and science, and read Harry Potter. Experimented with what you see in my comprehension scores this week."
$ws.Range("C12").Value = "Count of Monte Cristo (Text-only, French, New):28; 100 años de soledad (Text-only, Spanish, New):25; House of Cards (Audiovisual, English, Familiar):30; Dr. House (Audiovisual, English, New):31; MultiAnime PodCast 3x01 (Audio-only, Spanish, New):34; La Cotorrisa - Anecdotario 106 - No me querían en la iglesia. (Audio-only, Spanish, New):23; Club de cuervos (Audiovisual, Spanish, New):24; Somos. (Audiovisual, Spanish, New):21;"

# Update selection to C13 (the cell the author clicked next)
$ws.Range("C13").Select()
